$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "Previously added" sheet: append the 4 rows that used to live on the
#    "New" sheet (rows 187-190), carrying over their original values/links.
# ---------------------------------------------------------------------------
$wsPrev = $wb.Worksheets.Item("Previously added")

$prevRows = @(
    @{ A = "https://www.ss.com/msg/lv/real-estate/wood/daugavpils-and-reg/malinovas-pag/mknjc.html"; B = "18 500 €";  C = "Daugavpils un raj."; D = "3 ha.";  E = "44700010164"; F = 45946.60208333333 },
    @{ A = "https://www.ss.com/msg/lv/real-estate/wood/dobele-and-reg/berzes-pag/kkjpg.html";        B = "1 400 €";   C = "Dobele un raj.";     D = "1 ha.";  E = "46520030087"; F = 45946.40972222222 },
    @{ A = "https://www.ss.com/msg/lv/real-estate/wood/jelgava-and-reg/kalnciems/cfghfd.html";       B = "89 500 €";  C = "Jelgava un raj.";    D = "18 ha."; E = "54310030137"; F = 45946.49097222222 },
    @{ A = "https://www.ss.com/msg/lv/real-estate/wood/valka-and-reg/blomes-pag/lkijb.html";         B = "123 456 €"; C = "Valka un raj.";      D = "6 ha.";  E = "94460010165"; F = 45945.875 }
)

$startRow = 187
for ($i = 0; $i -lt $prevRows.Count; $i++) {
    $r = $startRow + $i
    $row = $prevRows[$i]

    # Copy the formatting from the last existing data row so the new rows
    # pick up the same styles (hyperlink font on A, normal font on B-E, date
    # format on F) as the rest of the sheet.
    $wsPrev.Range("A186:F186").Copy()
    $wsPrev.Range("A$r`:F$r").PasteSpecial(-4122)

    $wsPrev.Range("A$r").Value = $row.A
    $wsPrev.Range("B$r").Value = $row.B
    $wsPrev.Range("C$r").Value = $row.C
    $wsPrev.Range("D$r").Value = $row.D

    # Cadastre numbers are plain digit strings - force text so they are not
    # reinterpreted as numbers.
    $wsPrev.Range("E$r").NumberFormat = "@"
    $wsPrev.Range("E$r").Value = $row.E

    $wsPrev.Range("F$r").Value = $row.F

    $wsPrev.Hyperlinks.Add($wsPrev.Range("A$r"), $row.A)

    # Re-apply the formatting once more - Hyperlinks.Add stamps its own
    # style on the anchor cell, so restore the plain column-A style.
    $wsPrev.Range("A186:F186").Copy()
    $wsPrev.Range("A$r`:F$r").PasteSpecial(-4122)
}

# ---------------------------------------------------------------------------
# 2) "New" sheet: drop the old 4 rows of data and replace them with the
#    freshly scraped listings.
# ---------------------------------------------------------------------------
$wsNew = $wb.Worksheets.Item("New")

# Remove every existing hyperlink (and its relationship) on this sheet so we
# don't leave stale links pointing at the old URLs.
$wsNew.Range("A2").Hyperlinks.Delete()

$newRows = @(
    @{ A = "https://www.ss.com/msg/lv/real-estate/wood/madona-and-reg/sarkanu-pag/iojbh.html";      B = "20 000 €"; C = "Madona un raj."; D = "9 ha."; E = "70900070024"; F = 45947.54236111111 },
    @{ A = "https://www.ss.com/msg/lv/real-estate/wood/madona-and-reg/berzaunes-pag/niefh.html";     B = "19 871 €"; C = "Madona un raj."; D = "1 ha."; E = "70460020031"; F = 45946.975 },
    @{ A = "https://www.ss.com/msg/lv/real-estate/wood/tukums-and-reg/lapmezciema-nov/ocdgb.html";   B = "16 000 €"; C = "Tukums un raj."; D = "4 ha."; E = "906600497";  F = 45947.552777777775 },
    @{ A = "https://www.ss.com/msg/lv/real-estate/wood/tukums-and-reg/lestenes-pag/bximhb.html";     B = "18 000 €"; C = "Tukums un raj."; D = "5 ha."; E = "";           F = 45946.759722222225 }
)

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = 2 + $i
    $row = $newRows[$i]

    $wsNew.Range("A$r").Value = $row.A
    $wsNew.Range("B$r").Value = $row.B
    $wsNew.Range("C$r").Value = $row.C
    $wsNew.Range("D$r").Value = $row.D

    # Cadastre numbers are plain digit strings - force text so they are not
    # reinterpreted as numbers.
    $wsNew.Range("E$r").NumberFormat = "@"
    $wsNew.Range("E$r").Value = $row.E

    $wsNew.Range("F$r").Value = $row.F

    $wsNew.Hyperlinks.Add($wsNew.Range("A$r"), $row.A)

    # Re-apply the original per-column styling - Hyperlinks.Add and the
    # text NumberFormat tweak both stamp their own style on the cells they
    # touch, so pull the canonical styles back from the "Previously added"
    # sheet (column A = hyperlink style, B-E = plain text, F = date).
    $wsPrev.Range("A186:F186").Copy()
    $wsNew.Range("A$r`:F$r").PasteSpecial(-4122)
}
